$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 26; A = "Globo"; B = "Bom Dia Inter"; C = "Social"; D = "2025-04-03T08:03"; E = "Negativo"; F = "Mais um dia de fila em Campos. Beneficiários dormem na fila para atualização do CadÚnico. Repórter *ao vivo*. Informações sobre mutirão. Entrevistas com as pessoas na fila e muitas reclamações. Pessoas estão se sentindo humilhadas. " },
    @{ r = 27; A = "Globo"; B = "RJ TV 1"; C = "Social"; D = "2025-04-03T11:52"; E = "Negativo"; F = "Desorganização no mutirão em Campos. Beneficiários do CadÚnico reclamam de filas gigantes e da falta de estrutura no atendimento. Repórter *ao vivo* do ginásio da Fundação Municipal de Esportes, onde acontece segundo dia do mutirão. Imagens do ginásio com pessoas sendo atendidas e outras na arquibancada aguardando atendimento. Lá fora, muitas pessoas na longa fila. Equipe do Cadùnco chegou mais cedo e já distribuiu 800 senhas. Depoimentos de beneficiários com reclamações. Uma senhora diz que foi no Cras de Travessão e foi bem atendida mas, com o braço machucado, teve que enfrentar longa fila, mas foi bem atendida no ginásio. Uma entrevista que disse ser `"moradora da rua`" também reclamou e disse que as autoridades alegam que a cidade está um `"petáculo`". Entrevista com assessor Ruan Barros. Também foram distribuídas pulseiras separando os grupos prioritários. Assessor também disse que tem água e banheiro. *com nota* da prefeitura. Apresentadores também criticaram a organização da prefeitura. " }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.A
    $ws.Cells.Item($row.r, 2).Value = $row.B
    $ws.Cells.Item($row.r, 3).Value = $row.C
    $ws.Cells.Item($row.r, 4).Value = $row.D
    $ws.Cells.Item($row.r, 5).Value = $row.E
    $ws.Cells.Item($row.r, 6).Value = $row.F
}
